$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The product management data row (row 2) is being cleared out - the
# specific gearbox/valve attributes that used to live under the header
# row are removed while the header row (row 1) and the cell formatting
# stay intact.
$ws.Range("A2:N2").ClearContents()

# Mirror the user's workflow: select the full data row (as Excel does
# when you click the row header before deleting its contents) and let
# the row shrink back down to a normal single-line height now that the
# wrapped text is gone.
$ws.Rows.Item(2).RowHeight = 15
$ws.Range("A2:XFD2").Select()
